$wb = $excel.ActiveWorkbook

# --- "parameters" sheet: add reference_node parameter row ---
# NOTE: do this before touching/activating the "edge" sheet so that the
# final active tab ends up being "edge" (activeTab=1), matching the target.
$wsParams = $wb.Worksheets.Item("parameters")
$wsParams.Range("B2").Value = 1
$wsParams.Range("A5").Value = "reference_node"
$wsParams.Range("B5").Value = "node1"
$wsParams.Range("C5").Value = "reference for voltage angles"
$wsParams.Range("C6").Select()

# --- "edge" sheet: add resistance/reactance for the row 5 branch, and make it the active tab ---
$wsEdge = $wb.Worksheets.Item("edge")
$wsEdge.Range("F5").Value = 0.001
$wsEdge.Range("G5").Value = 0.02
$wsEdge.Activate()
$wsEdge.Range("D5").Select()
